# GSC export data update: the 2025-11-17 entry has aged out of the report
# window, so its row is removed from the "Chart" sheet and every
# subsequent row shifts up by one (dates stay aligned, the "Items" counts
# slide up to the following day's row).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 holds the oldest date (2025-11-17); deleting the whole row shifts
# all the rows below it up by one, which is exactly what the export diff
# shows (C2 becomes the old C3 value, etc.), and drops the now-unused
# "2025-11-17" shared string automatically.
$ws.Rows.Item(2).Delete()
